$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the E2/F2 probability values (0.5/0.6 -> 0.6/0.5)
$ws.Range("E2").Value = 0.6
$ws.Range("F2").Value = 0.5

# Clear the stray "applyNumberFormat" style that was left on E4:F4
# (value is unchanged, only the cell format goes back to the default)
$ws.Range("E4:F4").Style = "Normal"

# Match the default column width used when the sheet was re-saved
$ws.StandardWidth = 8.7265625

# Leave the active selection on F3, as in the edited workbook
$ws.Range("F3").Select()
